# Update "想去人数" (number of people interested) figures on the
# 展览 and 全部类型 sheets to the refreshed values from the data source.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value = 2136
    $ws.Range("F4").Value = 299
    $ws.Range("F6").Value = 6389
    $ws.Range("F7").Value = 270
}
